$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.014.51"
$ws.Range("E2").Value = "  -3.40%  "

$ws.Range("D3").Value = "2.556.99"
$ws.Range("E3").Value = "  -3.48%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'300.44"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").Value = "'94.41"
$ws.Range("E6").Value = "  -2.45%  "

$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  -2.16%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("D10").Value = "'36.35"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("D12").Value = "'7.75"
$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "'0.115"
$ws.Range("E13").Value = "  +8.13%  "

$ws.Range("D14").Value = "2.955.93"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").Value = "2.571.99"
$ws.Range("E15").Value = "  -3.02%  "

$ws.Range("D16").Value = "'0.882"
$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("D17").Value = "'14.22"
$ws.Range("E17").Value = "  -4.00%  "

$ws.Range("D18").Value = "43.054.40"
$ws.Range("E18").Value = "  -3.32%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0984"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.74"
$ws.Range("E20").Value = "  +2.28%  "

$ws.Range("D21").Value = "'6.56"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").Value = "'72.08"
$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("D23").Value = "'254.73"
$ws.Range("E23").Value = "  -7.19%  "

$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").Value = "'2.13"
$ws.Range("E25").Value = "  -5.97%  "

$ws.Range("D26").Value = "'28.91"
$ws.Range("E26").Value = "  -5.07%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'10.24"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").Value = "'36.96"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  -5.81%  "

$ws.Range("D31").Value = "'6.09"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").Value = "'152.72"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.77"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'3.38"
$ws.Range("E34").Value = "  -8.79%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'2.14"
$ws.Range("E35").Value = "  -6.69%  "

$ws.Range("D36").Value = "'0.0797"
$ws.Range("E36").Value = "  -3.37%  "

$ws.Range("D37").Value = "'0.114"
$ws.Range("E37").Value = "  -2.86%  "

$ws.Range("D38").Value = "'17.72"
$ws.Range("E38").Value = "  +12.36%  "

$ws.Range("D39").Value = "'0.119"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").Value = "'23.09"
$ws.Range("E40").Value = "  -7.58%  "

$ws.Range("D41").Value = "'2.30"
$ws.Range("E41").Value = "  +44.26%  "

$ws.Range("D42").Value = "'3.43"
$ws.Range("E42").Value = "  -3.14%  "

$ws.Range("D43").Value = "'0.0312"
$ws.Range("E43").Value = "  -2.17%  "

$ws.Range("D44").Value = "'3.87"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").Value = "2.105.93"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "'9.21"
$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").Value = "'84.76"
$ws.Range("E48").Value = "  -7.35%  "

$ws.Range("D49").Value = "2.810.99"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("D50").Value = "'105.59"
$ws.Range("E50").Value = "  -3.31%  "

$ws.Range("D51").Value = "'1.69"
$ws.Range("E51").Value = "  +0.09%  "
